$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Widen column D slightly (11 -> 12)
$ws.Columns.Item(4).ColumnWidth = 11.1666666666667

# Row 2 (OTROS): VENTA becomes -10.44, POR CUMPLIR becomes 10.44
$ws.Range("D2").Value = -10.44
$ws.Range("E2").Value = 10.44

# Row 4 (TOTAL): VENTA becomes -10.44, POR CUMPLIR becomes 20010.44, CUMPLIMIENTO becomes -0.000522
$ws.Range("D4").Value = -10.44
$ws.Range("E4").Value = 20010.44
$ws.Range("F4").Value = -0.000522
